$p = $ppt.ActivePresentation
try { Write-Host "DocumentTheme:" $p.DocumentTheme } catch { Write-Host "err:" $_ }
try { Write-Host "CustomXMLParts:" $p.CustomXMLParts } catch { Write-Host "err:" $_ }
